$d = $word.ActiveDocument

# 1. Insert " (Kim)" right after "Vergleich von Darcy" in the Darcy/Cozeny-Karman
#    paragraph ("- Vergleich von Darcy mit bspw. Cozeny-Karman ...").
$find = $d.Content
$find.Find.Execute("Vergleich von Darcy") | Out-Null
$darcyParaIndex = $find.Paragraphs(1).Index
$find.Collapse(0)  # wdCollapseEnd
$find.InsertAfter(" (Kim)")

# 2. The _GoBack bookmark moves from the title paragraph ("Transportprozesse") to a
#    brand new paragraph ("Cozeny Karman (Ecequel)") inserted right after the Darcy
#    paragraph (i.e. right before "3. Numerische Modelle").
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$darcyPara = $d.Paragraphs($darcyParaIndex)
$darcyPara.Range.InsertParagraphAfter()
$newParaIndex = $darcyParaIndex + 1
$newPara = $d.Paragraphs($newParaIndex).Range
$newPara.Text = "Cozeny Karman (Ecequel)"

$newParaStart = $d.Paragraphs($newParaIndex).Range
$newParaStart.Collapse(1)  # wdCollapseStart
$d.Bookmarks.Add("_GoBack", $newParaStart) | Out-Null
